$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the paragraph that currently holds the _GoBack bookmark
# ("I have began implementing feature.js ...") and the two empty
# trailing paragraphs that follow it (right before the sectPr).
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)       # second (last) empty paragraph
$beforeLastPara = $d.Paragraphs.Item($count - 1)  # first empty paragraph - stays as-is (matches target "<w:p/>")

# Helper: always re-fetch the still-empty "last paragraph" (its index
# never changes because we only ever insert new paragraphs *before* it).
function Get-LastParaRange {
    $p = $d.Paragraphs.Item($d.Paragraphs.Count)
    return $p.Range
}

# ------------------------------------------------------------------
# Paragraph: "<empty>" (plain <w:p/>)
# ------------------------------------------------------------------
$r = Get-LastParaRange
$r.Collapse(1)
[void]$r.InsertParagraphBefore()

# ------------------------------------------------------------------
# Paragraph: "15th February 2016" (bold, "th" superscript, bold
# paragraph mark)
# ------------------------------------------------------------------
$r = Get-LastParaRange
$r.Collapse(1)
[void]$r.InsertParagraphBefore()
$p = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$pr = $p.Range
$paraStart = $pr.Start
$pr.Text = "15th February 2016"
$d.Range($paraStart, $paraStart + 19).Font.Bold = $true
$d.Range($paraStart + 2, $paraStart + 4).Font.Superscript = $true
$markRange = $p.Range
$markRange.Collapse(0)
$markRange.Font.Bold = $true

# ------------------------------------------------------------------
# Paragraph: empty, but with a bold paragraph mark
# ------------------------------------------------------------------
$r = Get-LastParaRange
$r.Collapse(1)
[void]$r.InsertParagraphBefore()
$p = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$markRange = $p.Range
$markRange.Collapse(0)
$markRange.Font.Bold = $true

# ------------------------------------------------------------------
# Paragraph: "I have began developing my website. I have started
# with the homepage." + two manual line breaks + "16th February
# 2016" (bold, "th" superscript)
# ------------------------------------------------------------------
$r = Get-LastParaRange
$r.Collapse(1)
[void]$r.InsertParagraphBefore()
$p = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$pr = $p.Range
$pr.Text = "I have began developing my website. I have started with the homepage."

$insPos = $p.Range.End - 1
$d.Range($insPos, $insPos).InsertBreak(6)

$insPos = $p.Range.End - 1
$d.Range($insPos, $insPos).InsertBreak(6)

$insPos = $p.Range.End - 1
$dateStart = $insPos
$d.Range($insPos, $insPos).InsertAfter("16th February 2016")

$d.Range($dateStart, $dateStart + 19).Font.Bold = $true
$d.Range($dateStart + 2, $dateStart + 4).Font.Superscript = $true

# ------------------------------------------------------------------
# Paragraph: "<empty>" (plain <w:p/>)
# ------------------------------------------------------------------
$r = Get-LastParaRange
$r.Collapse(1)
[void]$r.InsertParagraphBefore()

# ------------------------------------------------------------------
# Final paragraph: "I have began developing the services page for
# my company website. " -- this re-uses the document's permanent
# trailing paragraph mark, and also receives the relocated
# _GoBack bookmark (Word only keeps a single _GoBack, so adding it
# here automatically removes it from its old location).
# ------------------------------------------------------------------
$lastP = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastP.Range.Text = "I have began developing the services page for my company website. "
$d.Bookmarks.Add("_GoBack", $lastP.Range)
